$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update header row (A1:D1) with new machine-friendly column names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Normalize capitalization of Spanish connector words (de, del, la, las,
#    los, el, y) inside state/municipality names -- e.g. "Pabellón de Arteaga"
#    becomes "Pabellón De Arteaga".
$renames = @(
    @{Cell='B7'; Value='Pabellón De Arteaga'},
    @{Cell='B8'; Value='Rincón De Romos'},
    @{Cell='B9'; Value='San Francisco De Los Romo'},
    @{Cell='B22'; Value='Amatenango De La Frontera'},
    @{Cell='B23'; Value='Amatenango Del Valle'},
    @{Cell='B26'; Value='Bejucal De Ocampo'},
    @{Cell='B32'; Value='Chiapa De Corzo'},
    @{Cell='B37'; Value='Comitán De Domínguez'},
    @{Cell='B59'; Value='Mazapa De Madero'},
    @{Cell='B63'; Value='Ocozocoautla De Espinosa'},
    @{Cell='B72'; Value='Salto De Agua'},
    @{Cell='B73'; Value='San Cristóbal De Las Casas'},
    @{Cell='B101'; Value='Hidalgo Del Parral'},
    @{Cell='B109'; Value='San Francisco De Conchos'},
    @{Cell='B124'; Value='San Juan De Sabinas'},
    @{Cell='A132'; Value='Ciudad De México'},
    @{Cell='B136'; Value='Cuajimalpa De Morelos'},
    @{Cell='B157'; Value='Nombre De Dios'},
    @{Cell='B165'; Value='San Juan De Guadalupe'},
    @{Cell='A172'; Value='Estado De México'},
    @{Cell='B172'; Value='Acambay De Ruíz Castañeda'},
    @{Cell='B174'; Value='Almoloya De Alquisiras'},
    @{Cell='B175'; Value='Almoloya De Juárez'},
    @{Cell='B179'; Value='Atizapán De Zaragoza'},
    @{Cell='B189'; Value='Ecatepec De Morelos'},
    @{Cell='B193'; Value='Ixtapan De La Sal'},
    @{Cell='B194'; Value='Ixtapan Del Oro'},
    @{Cell='B201'; Value='Naucalpan De Juárez'},
    @{Cell='B207'; Value='San Felipe Del Progreso'},
    @{Cell='B208'; Value='San Martín De Las Pirámides'},
    @{Cell='B215'; Value='Tenango Del Valle'},
    @{Cell='B222'; Value='Tlalnepantla De Baz'},
    @{Cell='B227'; Value='Valle De Chalco Solidaridad'},
    @{Cell='B236'; Value='San Miguel De Allende'},
    @{Cell='B237'; Value='Apaseo El Alto'},
    @{Cell='B238'; Value='Apaseo El Grande'},
    @{Cell='B244'; Value='Dolores Hidalgo Cuna De La Independencia Nacional'},
    @{Cell='B248'; Value='Jaral Del Progreso'},
    @{Cell='B258'; Value='San Diego De La Unión'},
    @{Cell='B260'; Value='San Francisco Del Rincón'},
    @{Cell='B262'; Value='San Luis De La Paz'},
    @{Cell='B263'; Value='Santa Cruz De Juventino Rosas'},
    @{Cell='B268'; Value='Valle De Santiago'},
    @{Cell='B274'; Value='Acapulco De Juárez'},
    @{Cell='B276'; Value='Ajuchitlán Del Progreso'},
    @{Cell='B277'; Value='Alcozauca De Guerrero'},
    @{Cell='B280'; Value='Atenango Del Río'},
    @{Cell='B281'; Value='Atlamajalcingo Del Monte'},
    @{Cell='B283'; Value='Atoyac De Álvarez'},
    @{Cell='B284'; Value='Ayutla De Los Libres'},
    @{Cell='B287'; Value='Buenavista De Cuéllar'},
    @{Cell='B288'; Value='Chilapa De Álvarez'},
    @{Cell='B289'; Value='Chilpancingo De Los Bravo'},
    @{Cell='B292'; Value='Coyuca De Benítez'},
    @{Cell='B293'; Value='Coyuca De Catalán'},
    @{Cell='B296'; Value='Cuetzala Del Progreso'},
    @{Cell='B297'; Value='Cutzamala De Pinzón'},
    @{Cell='B302'; Value='Huitzuco De Los Figueroa'},
    @{Cell='B303'; Value='Iguala De La Independencia'},
    @{Cell='B307'; Value='La Unión De Isidoro Montes De Oca'},
    @{Cell='B320'; Value='Taxco De Alarcón'},
    @{Cell='B322'; Value='Técpan De Galeana'},
    @{Cell='B324'; Value='Tepecoacuilco De Trujano'},
    @{Cell='B327'; Value='Tlapa De Comonfort'},
    @{Cell='B336'; Value='Agua Blanca De Iturbide'},
    @{Cell='B340'; Value='Atotonilco El Grande'},
    @{Cell='B345'; Value='Cuautepec De Hinojosa'},
    @{Cell='B349'; Value='Huasca De Ocampo'},
    @{Cell='B351'; Value='Huejutla De Reyes'},
    @{Cell='B354'; Value='Jacala De Ledezma'},
    @{Cell='B359'; Value='Mixquiahuala De Juárez'},
    @{Cell='B360'; Value='Molango De Escamilla'},
    @{Cell='B362'; Value='Omitlán De Juárez'},
    @{Cell='B363'; Value='Pachuca De Soto'},
    @{Cell='B365'; Value='Progreso De Obregón'},
    @{Cell='B372'; Value='Tenango De Doria'},
    @{Cell='B374'; Value='Tepehuacán De Guerrero'},
    @{Cell='B375'; Value='Tepeji Del Río De Ocampo'},
    @{Cell='B382'; Value='Tula De Allende'},
    @{Cell='B383'; Value='Tulancingo De Bravo'},
    @{Cell='B386'; Value='Zapotlán De Juárez'},
    @{Cell='B390'; Value='Ahualulco De Mercado'},
    @{Cell='B393'; Value='Atemajac De Brizuela'},
    @{Cell='B395'; Value='Atotonilco El Alto'},
    @{Cell='B396'; Value='Autlán De Navarro'},
    @{Cell='B403'; Value='Encarnación De Díaz'},
    @{Cell='B405'; Value='Huejuquilla El Alto'},
    @{Cell='B409'; Value='Lagos De Moreno'},
    @{Cell='B412'; Value='Ojuelos De Jalisco'},
    @{Cell='B415'; Value='San Juan De Los Lagos'},
    @{Cell='B416'; Value='San Martín De Bolaños'},
    @{Cell='B418'; Value='San Miguel El Alto'},
    @{Cell='B419'; Value='San Sebastián Del Oeste'},
    @{Cell='B420'; Value='Santa María Del Oro'},
    @{Cell='B421'; Value='Tamazula De Gordiano'},
    @{Cell='B423'; Value='Tepatitlán De Morelos'},
    @{Cell='B424'; Value='Tizapán El Alto'},
    @{Cell='B425'; Value='Tlajomulco De Zúñiga'},
    @{Cell='B432'; Value='Unión De Tula'},
    @{Cell='B434'; Value='Yahualica De González Gallo'},
    @{Cell='B437'; Value='Zapotlán El Grande'},
    @{Cell='B453'; Value='Cojumatlán De Régules'},
    @{Cell='B511'; Value='Coatlán Del Río'},
    @{Cell='B519'; Value='Puente De Ixtla'},
    @{Cell='B521'; Value='Tetela Del Volcán'},
    @{Cell='B522'; Value='Tlaltizapán De Zapata'},
    @{Cell='B532'; Value='Ixtlán Del Río'},
    @{Cell='B537'; Value='Santa María Del Oro'},
    @{Cell='B549'; Value='Mier Y Noriega'},
    @{Cell='B553'; Value='San Nicolás De Los Garza'},
    @{Cell='B555'; Value='Acatlán De Pérez Figueroa'},
    @{Cell='B561'; Value='Chalcatongo De Hidalgo'},
    @{Cell='B562'; Value='Coicoyán De Las Flores'},
    @{Cell='B563'; Value='Constancia Del Rosario'},
    @{Cell='B566'; Value='Heroica Ciudad De Huajuapan De León'},
    @{Cell='B567'; Value='Heroica Ciudad De Tlaxiaco'},
    @{Cell='B568'; Value='Ixtlán De Juárez'},
    @{Cell='B569'; Value='Heroica Ciudad De Juchitán De Zaragoza'},
    @{Cell='B572'; Value='Mariscala De Juárez'},
    @{Cell='B574'; Value='Miahuatlán De Porfirio Díaz'},
    @{Cell='B575'; Value='Oaxaca De Juárez'},
    @{Cell='B576'; Value='Ocotlán De Morelos'},
    @{Cell='B578'; Value='Putla Villa De Guerrero'},
    @{Cell='B581'; Value='San Antonino El Alto'},
    @{Cell='B590'; Value='San Francisco Del Mar'},
    @{Cell='B602'; Value='San Juan Bautista Lo De Soto'},
    @{Cell='B619'; Value='San Miguel Del Puerto'},
    @{Cell='B628'; Value='San Pablo Villa De Mitla'},
    @{Cell='B631'; Value='San Pedro El Alto'},
    @{Cell='B653'; Value='Santa Cruz Tacache De Mina'},
    @{Cell='B657'; Value='Santa Inés De Zaragoza'},
    @{Cell='B686'; Value='Santo Domingo De Morelos'},
    @{Cell='B692'; Value='Tataltepec De Valdés'},
    @{Cell='B693'; Value='Teotitlán De Flores Magón'},
    @{Cell='B694'; Value='Teotitlán Del Valle'},
    @{Cell='B696'; Value='Tlacolula De Matamoros'},
    @{Cell='B697'; Value='Villa De Etla'},
    @{Cell='B698'; Value='Villa De Tututepec De Melchor Ocampo'},
    @{Cell='B699'; Value='Villa De Zaachila'},
    @{Cell='B700'; Value='Villa Sola De Vega'},
    @{Cell='B703'; Value='Zimatlán De Álvarez'},
    @{Cell='B717'; Value='Chalchicomula De Sesma'},
    @{Cell='B721'; Value='Chila De La Sal'},
    @{Cell='B725'; Value='Cuayuca De Andrade'},
    @{Cell='B726'; Value='Cuetzalan Del Progreso'},
    @{Cell='B732'; Value='Huehuetlán El Chico'},
    @{Cell='B733'; Value='Huehuetlán El Grande'},
    @{Cell='B736'; Value='Ixcamilpa De Guerrero'},
    @{Cell='B737'; Value='Izúcar De Matamoros'},
    @{Cell='B741'; Value='Los Reyes De Juárez'},
    @{Cell='B745'; Value='Palmar De Bravo'},
    @{Cell='B757'; Value='Tecali De Herrera'},
    @{Cell='B761'; Value='Tepanco De López'},
    @{Cell='B763'; Value='Tepexi De Rodríguez'},
    @{Cell='B764'; Value='Tetela De Ocampo'},
    @{Cell='B769'; Value='Tlacotepec De Benito Juárez'},
    @{Cell='B781'; Value='Amealco De Bonfil'},
    @{Cell='B783'; Value='Cadereyta De Montes'},
    @{Cell='B786'; Value='Jalpan De Serra'},
    @{Cell='B787'; Value='Landa De Matamoros'},
    @{Cell='B790'; Value='San Juan Del Río'},
    @{Cell='B804'; Value='Ciudad Del Maíz'},
    @{Cell='B815'; Value='San Ciro De Acosta'},
    @{Cell='B819'; Value='Santa María Del Río'},
    @{Cell='B821'; Value='Soledad De Graciano Sánchez'},
    @{Cell='B826'; Value='Tanquián De Escobedo'},
    @{Cell='B828'; Value='Villa De Guadalupe'},
    @{Cell='B829'; Value='Villa De La Paz'},
    @{Cell='B830'; Value='Villa De Ramos'},
    @{Cell='B880'; Value='Soto La Marina'},
    @{Cell='B893'; Value='Ixtacuixtla De Mariano Matamoros'},
    @{Cell='B896'; Value='Sanctórum De Lázaro Cárdenas'},
    @{Cell='B898'; Value='Tetla De La Solidaridad'},
    @{Cell='B913'; Value='Amatlán De Los Reyes'},
    @{Cell='B923'; Value='Castillo De Teayo'},
    @{Cell='B925'; Value='Cazones De Herrera'},
    @{Cell='B939'; Value='Cosamaloapan De Carpio'},
    @{Cell='B940'; Value='Cosautlán De Carvajal'},
    @{Cell='B953'; Value='Hueyapan De Ocampo'},
    @{Cell='B954'; Value='Ignacio De La Llave'},
    @{Cell='B956'; Value='Ixhuatlán Del Café'},
    @{Cell='B963'; Value='Juchique De Ferrer'},
    @{Cell='B966'; Value='Las Vigas De Ramírez'},
    @{Cell='B967'; Value='Lerdo De Tejada'},
    @{Cell='B972'; Value='Martínez De La Torre'},
    @{Cell='B977'; Value='Nanchital De Lázaro Cárdenas Del Río'},
    @{Cell='B986'; Value='Paso De Ovejas'},
    @{Cell='B987'; Value='Paso Del Macho'},
    @{Cell='B990'; Value='Poza Rica De Hidalgo'},
    @{Cell='B997'; Value='Sayula De Alemán'},
    @{Cell='B999'; Value='Soledad De Doblado'},
    @{Cell='B1029'; Value='Vega De Alatorre'},
    @{Cell='B1037'; Value='Zontecomatlán De López Y Fuentes'},
    @{Cell='B1041'; Value='Cañitas De Felipe Pescador'},
    @{Cell='B1053'; Value='Mezquital Del Oro'},
    @{Cell='B1055'; Value='Moyahua De Estrada'},
    @{Cell='B1056'; Value='Noria De Ángeles'}
)

foreach ($item in $renames) {
    $ws.Range($item.Cell).Value = $item.Value
}

# 3. Remove the trailing metadata/footnote rows (1072-1076), which shifts
#    the sheet's used range/dimension down to row 1070.
$ws.Rows("1072:1076").Delete()
